$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.791.96"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +1.90%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.860.40"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.65%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9995"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "244.15"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.40%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6432"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +4.60%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07539"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +2.59%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2988"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +1.92%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.62"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +5.83%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07696"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.68%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.873.91"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +2.33%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.057"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +1.37%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6924"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +2.68%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "84.10"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +1.82%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.000009866"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +10.45%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.124"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +4.23%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "29.795.99"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +1.93%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.113.76"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +1.11%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "236.41"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.16%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.66"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +1.12%  "
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.04%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.578"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +2.58%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.08%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "158.45"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.16%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1424"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +2.28%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.570"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.20%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.95"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +1.71%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.06207"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +7.08%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.493"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.18%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.289"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.159"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +1.32%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.097"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +0.20%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.898"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +2.21%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.175"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +3.40%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7324"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +1.55%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.610"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.14%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.824"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -1.38%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01793"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +1.75%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.220.07"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.23%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.322"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +1.83%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9212"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +1.50%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.001"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.03%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.021.80"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.75%  "
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.18%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "67.12"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +2.14%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000123"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +2.89%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4072"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.82%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.170"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.16%  "
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +5.34%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.1131"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -3.77%  "
